# Auto-generated edit script: update countries list order + refresh covid stats + timestamp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 13:28"

# Full country data table (rows 4.. ), columns A.."H (Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$data = @(
  ,("Estados Unidos",5874295,149,3167164,2526526,0,1,180605)
  ,("Brasil",3605783,0,2709638,781373,0,0,114772)
  ,("India",3110761,5576,2338899,714147,0,23,57715)
  ,("Rusia",961493,4744,773095,171950,0,65,16448)
  ,("Sudafrica",609773,0,506470,90244,0,0,13059)
  ,("Peru",594326,0,399357,167306,0,0,27663)
  ,("Mexico",560164,3948,383872,115812,0,226,60480)
  ,("Colombia",541147,0,374030,149801,0,0,17316)
  ,("España",407879,0,0,0,0,0,28838)
  ,("Chile",397665,0,371179,15634,0,0,10852)
  ,("Iran",361150,2245,311365,29009,0,133,20776)
  ,("Argentina",342154,0,251400,83769,0,0,6985)
  ,("Reino Unido",325642,0,0,0,0,0,41429)
  ,("Arabia Saudita",307479,0,280143,23687,0,0,3649)
  ,("Banglades",297083,2485,182875,110225,0,42,3983)
  ,("Pakistan",293261,496,276829,10188,0,9,6244)
  ,("Italia",259345,0,205470,18438,0,0,35437)
  ,("Turquia",258249,0,237165,14963,0,0,6121)
  ,("Francia",242899,0,84973,127413,0,0,30513)
  ,("Alemania",234637,148,209600,15705,0,0,9332)
  ,("Irak",204341,0,146409,51504,0,0,6428)
  ,("Filipinas",194252,4686,132042,59200,0,13,3010)
  ,("Indonesia",155412,1877,111060,37593,0,79,6759)
  ,("Canada",124896,0,111112,4711,0,0,9073)
  ,("Catar",117266,258,114099,2973,0,1,194)
  ,("Bolivia",109149,722,45396,59244,0,67,4509)
  ,("Ecuador",107769,0,94816,6643,0,0,6310)
  ,("Ucrania",106754,1799,54524,49937,0,22,2293)
  ,("Kazajistan",104718,175,92598,10705,0,0,1415)
  ,("Israel",103274,611,80521,21914,0,5,839)
  ,("Egipto",97340,0,65927,26151,0,0,5262)
  ,("Republica Dominicana",91161,0,60670,28924,0,0,1567)
  ,("Panama",86900,0,61420,23588,0,0,1892)
  ,("Suecia",86068,0,0,0,0,0,5810)
  ,("China",84967,16,79925,408,0,0,4634)
  ,("Oman",84509,740,78912,4960,0,28,637)
  ,("Belgica",81936,468,18225,53719,0,4,9992)
  ,("Kuwait",80528,0,72307,7706,0,0,515)
  ,("Rumania",79330,825,35517,40504,0,37,3309)
  ,("Bielorrusia",70645,177,68925,1074,0,4,646)
  ,("Guatemala",68188,0,56778,8816,0,0,2594)
  ,("Emiratos Arabes Unidos",67282,275,58582,8324,0,1,376)
  ,("Paises Bajos",66554,0,0,0,0,0,6200)
  ,("Polonia",62310,548,42448,17902,0,5,1960)
  ,("Japon",61747,0,48550,12021,0,0,1176)
  ,("Singapur",56404,51,54164,2213,0,0,27)
  ,("Portugal",55597,0,40774,13027,0,0,1796)
  ,("Honduras",54511,528,8532,44325,0,11,1654)
  ,("Marruecos",52349,0,36343,15118,0,0,888)
  ,("Nigeria",52227,0,38945,12280,0,0,1002)
  ,("Barein",49330,0,45981,3165,0,0,184)
  ,("Ghana",43505,0,41532,1712,0,0,261)
  ,("Kirguistan",43126,103,36615,5454,0,1,1057)
  ,("Armenia",42825,33,36049,5922,0,2,854)
  ,("Argelia",41460,0,29142,10883,0,0,1435)
  ,("Etiopia",40671,0,14995,24998,0,0,678)
  ,("Suiza",40060,157,34400,3659,0,0,2001)
  ,("Venezuela",39564,0,29966,9269,0,0,329)
  ,("Uzbekistan",39065,119,34987,3803,0,2,275)
  ,("Afganistan",38054,55,28360,8305,0,2,1389)
  ,("Azerbaiyan",35274,0,32993,1763,0,0,518)
  ,("Costa Rica",33820,0,10518,22947,0,0,355)
  ,("Moldavia",33478,0,22683,9855,0,0,940)
  ,("Nepal",32678,743,18806,13715,0,8,157)
  ,("Kenia",32364,0,18670,13146,0,0,548)
  ,("Serbia",30657,0,28275,1684,0,0,698)
  ,("Irlanda",27969,0,23364,2828,0,0,1777)
  ,("Austria",25495,242,21657,3105,0,1,733)
  ,("Australia",24916,104,19603,4796,0,15,517)
  ,("El Salvador",24811,189,12492,11650,0,8,669)
  ,("Chequia",21923,0,16125,5386,0,0,412)
  ,("Estado de Palestina",19213,411,11870,7213,0,2,130)
  ,("Camerun",18762,0,16540,1814,0,0,408)
  ,("Bosnia y Herzegovina",17715,0,11329,5854,0,0,532)
  ,("Corea del Sur",17665,266,14219,3137,0,0,309)
  ,("Costa de Marfil",17471,0,15301,2057,0,0,113)
  ,("Dinamarca",16317,0,14180,1515,0,0,622)
  ,("Bulgaria",15287,0,10338,4404,0,0,545)
  ,("Madagascar",14402,75,13436,788,0,0,178)
  ,("Republica de Macedonia",13595,0,10110,2921,0,0,564)
  ,("Paraguay",13233,0,7417,5611,0,0,205)
  ,("Senegal",13013,64,8595,4146,0,3,272)
  ,("Sudan",12836,0,6497,5524,0,0,815)
  ,("Libano",12698,0,3625,8950,0,0,123)
  ,("Zambia",11082,0,9942,860,0,0,280)
  ,("Libia",11009,572,1096,9714,0,11,199)
  ,("Noruega",10323,0,9150,909,0,0,264)
  ,("Consejo Danes para los Refugiados",9830,0,8934,645,0,0,251)
  ,("Malasia",9274,7,8965,184,0,0,125)
  ,("Guinea",9013,46,7823,1136,0,1,54)
  ,("Guayana Francesa",8797,0,8307,435,0,0,55)
  ,("Grecia",8664,0,3804,4618,0,0,242)
  ,("Albania",8427,0,4332,3845,0,0,250)
  ,("Gabon",8388,0,6734,1601,0,0,53)
  ,("Tayikistan",8311,0,7108,1137,0,0,66)
  ,("Croacia",8175,0,5801,2203,0,0,171)
  ,("Haiti",8082,0,5624,2262,0,0,196)
  ,("Finlandia",7938,18,7100,504,0,0,334)
  ,("Luxemburgo",7775,0,6969,682,0,0,124)
  ,("Mauritania",6905,0,6232,515,0,0,158)
  ,("Maldivas",6779,0,4222,2531,0,0,26)
  ,("Zimbabue",5930,0,4872,903,0,0,155)
  ,("Namibia",5854,0,2509,3293,0,0,52)
  ,("Malaui",5414,0,3012,2234,0,0,168)
  ,("Republica de Yibuti",5382,0,5254,68,0,0,60)
  ,("Hungria",5191,36,3695,883,0,0,613)
  ,("Guinea Ecuatorial",4926,0,3795,1048,0,0,83)
  ,("Hong Kong",4692,9,4052,563,0,0,77)
  ,("Republica de Africa Central",4679,0,1755,2863,0,0,61)
  ,("Montenegro",4343,0,3356,903,0,0,84)
  ,("Nicaragua",4311,0,2913,1265,0,0,133)
  ,("Suazilandia",4225,0,2898,1242,0,0,85)
  ,("Congo",3850,0,1625,2148,0,0,77)
  ,("Cuba",3682,0,3044,547,0,0,91)
  ,("Surinam",3607,0,2688,861,0,0,58)
  ,("Cabo Verde",3509,0,2540,932,0,0,37)
  ,("Eslovaquia",3424,68,2153,1238,0,0,33)
  ,("Tailandia",3397,2,3222,117,0,0,58)
  ,("Mozambique",3395,0,1503,1872,0,0,20)
  ,("Somalia",3269,0,2396,780,0,0,93)
  ,("Mayotte",3237,0,2964,234,0,0,39)
  ,("Ruanda",3089,0,1755,1322,0,0,12)
  ,("Sri Lanka",2953,0,2811,130,0,0,12)
  ,("Tunez",2818,0,1443,1304,0,0,71)
  ,("Mali",2705,0,2018,562,0,0,125)
  ,("Gambia",2685,0,490,2108,0,0,87)
  ,("Lituania",2673,38,1766,822,0,1,85)
  ,("Eslovenia",2665,14,2122,410,0,2,133)
  ,("Sudan del Sur",2499,0,1290,1162,0,0,47)
  ,("Estonia",2275,3,2025,186,0,1,64)
  ,("Uganda",2263,0,1199,1044,0,0,20)
  ,("Siria",2217,0,505,1623,0,0,89)
  ,("Angola",2171,0,818,1257,0,0,96)
  ,("Guinea-Bisau",2149,0,1015,1101,0,0,33)
  ,("Benin",2115,0,1705,371,0,0,39)
  ,("Islandia",2064,0,1939,115,0,0,10)
  ,("Sierra Leona",1992,0,1550,373,0,0,69)
  ,("Yemen",1911,0,1086,272,0,0,553)
  ,("Bahamas",1765,0,227,1509,0,0,29)
  ,("Nueva Zelanda",1683,9,1538,123,0,0,22)
  ,("Malta",1667,55,977,680,0,0,10)
  ,("Jordania",1609,0,1320,277,0,0,12)
  ,("Aruba",1568,0,455,1106,0,0,7)
  ,("Jamaica",1529,116,819,694,0,0,16)
  ,("Uruguay",1527,0,1276,209,0,0,42)
  ,("Republica de Chipre",1421,0,878,523,0,0,20)
  ,("Georgia",1421,10,1137,267,0,0,17)
  ,("Letonia",1337,0,1093,211,0,0,33)
  ,("Burkina Faso",1320,0,1043,222,0,0,55)
  ,("Botsuana",1308,0,136,1169,0,0,3)
  ,("Liberia",1286,0,816,388,0,0,82)
  ,("Togo",1277,0,910,340,0,0,27)
  ,("Reunion",1209,0,692,511,0,0,6)
  ,("Niger",1172,0,1083,20,0,0,69)
  ,("Principado de Andorra",1045,0,875,117,0,0,53)
  ,("Vietnam",1022,6,587,408,0,0,27)
  ,("Lesoto",1015,0,472,513,0,0,30)
  ,("Trinidad yTobago",1007,0,165,828,0,0,14)
  ,("Republica del Chad",986,0,870,40,0,0,76)
  ,("Guyana",955,0,490,434,0,0,31)
  ,("Santo Tome y Principe",892,0,831,46,0,0,15)
  ,("Guadalupe",771,0,289,467,0,0,15)
  ,("Crucero",712,0,651,48,0,0,13)
  ,("San Marino",704,0,657,5,0,0,42)
  ,("Belice",686,0,44,636,0,0,6)
  ,("Tanzania",509,0,183,305,0,0,21)
  ,("Taiwan",487,0,457,23,0,0,7)
  ,("Martinica",464,0,98,350,0,0,16)
  ,("Birmania",463,13,341,116,0,0,6)
  ,("Burundi",430,0,336,93,0,0,1)
  ,("Comoras",417,0,396,14,0,0,7)
  ,("Islas Feroe",410,26,337,73,0,0,0)
  ,("Papua Nueva Guinea",401,0,232,165,0,0,4)
  ,("San Martin (Parte Holandesa)",396,0,147,232,0,0,17)
  ,("Islas Turcas y Caicos",383,0,102,279,0,0,2)
  ,("Mauricio",346,0,335,1,0,0,10)
  ,("Isla de Man",336,0,312,0,0,0,24)
  ,("Eritrea",306,0,274,32,0,0,0)
  ,("Polinesia Francesa",298,0,148,150,0,0,0)
  ,("Mongolia",298,0,289,9,0,0,0)
  ,("Camboya",273,0,263,10,0,0,0)
  ,("Gibraltar",248,2,203,45,0,0,0)
  ,("Islas Caimanes",205,0,202,2,0,0,1)
  ,("San Martin (Parte Francesa)",176,0,52,119,0,0,5)
  ,("Bermudas",167,0,149,9,0,0,9)
  ,("Barbados",161,0,132,22,0,0,7)
  ,("Butan",155,0,115,40,0,0,0)
  ,("Brunei",143,0,139,1,0,0,3)
  ,("Seychelles",132,0,126,6,0,0,0)
  ,("Monaco",112,0,83,28,0,0,1)
  ,("Liechtenstein",100,1,90,9,0,0,1)
  ,("Antigua y Barbuda",94,0,89,2,0,0,3)
  ,("San Vicente y las Granadinas",58,0,56,2,0,0,0)
  ,("Macao",46,0,46,0,0,0,0)
  ,("Curazao",43,4,34,8,0,0,1)
  ,("Puerto Rico",39,0,1,36,0,0,2)
  ,("Guam",32,0,0,31,0,0,1)
  ,("Fiyi",28,0,23,4,0,0,1)
  ,("Santa Lucia",26,0,25,1,0,0,0)
  ,("Timor Oriental",26,0,25,1,0,0,0)
  ,("Granada",24,0,24,0,0,0,0)
  ,("Nueva Caledonia",23,0,23,0,0,0,0)
  ,("Laos",22,0,20,2,0,0,0)
  ,("Islas Virgenes Britanicas",21,0,8,12,0,0,1)
  ,("Dominica",19,0,18,1,0,0,0)
  ,("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0)
  ,("San Cristobal y Nieves",17,0,17,0,0,0,0)
  ,("San Bartolome",16,0,9,7,0,0,0)
  ,("Groenlandia",14,0,14,0,0,0,0)
  ,("Bonaire, San Eustaquio y Saba",13,0,7,6,0,0,0)
  ,("Montserrat",13,0,12,0,0,0,1)
  ,("Islas Malvinas",13,0,13,0,0,0,0)
  ,("Santa Sede",12,0,12,0,0,0,0)
  ,("Sahara Occidental",10,0,8,1,0,0,1)
  ,("San Pedro y Miquelon",5,0,1,4,0,0,0)
  ,("Anguila",3,0,3,0,0,0,0)
)

$r = 4
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $r = $r + 1
}

"Done: wrote $($data.Count) rows"